{"js": "// Replace the arithmetic expressions (three-digit x one-digit multiplication\n// problems) in the document's table with the new set of expressions, as\n// described by the commit diff. Each old expression is unique in the\n// document, so a plain text search/replace for each pair is unambiguous.\n\nconst replacements = [\n  [\"578\u00d74=2312\", \"344\u00d77=2408\"],\n  [\"368\u00d78=2944\", \"686\u00d76=4116\"],\n  [\"900\u00d76=5400\", \"986\u00d75=4930\"],\n  [\"804\u00d79=7236\", \"193\u00d72=386\"],\n  [\"766\u00d76=4596\", \"996\u00d72=1992\"],\n  [\"425\u00d73=1275\", \"576\u00d72=1152\"],\n  [\"326\u00d78=2608\", \"405\u00d78=3240\"],\n  [\"266\u00d79=2394\", \"463\u00d78=3704\"],\n  [\"810\u00d73=2430\", \"998\u00d78=7984\"],\n  [\"452\u00d76=2712\", \"460\u00d75=2300\"],\n  [\"346\u00d74=1384\", \"934\u00d76=5604\"],\n  [\"822\u00d76=4932\", \"199\u00d77=1393\"],\n  [\"715\u00d75=3575\", \"613\u00d78=4904\"],\n  [\"187\u00d73=561\", \"881\u00d78=7048\"],\n  [\"178\u00d77=1246\", \"667\u00d72=1334\"],\n  [\"305\u00d72=610\", \"391\u00d73=1173\"],\n  [\"581\u00d74=2324\", \"238\u00d72=476\"],\n  [\"244\u00d72=488\", \"803\u00d78=6424\"],\n  [\"705\u00d77=4935\", \"333\u00d74=1332\"],\n  [\"564\u00d73=1692\", \"967\u00d72=1934\"],\n  [\"583\u00d75=2915\", \"131\u00d76=786\"],\n  [\"357\u00d77=2499\", \"442\u00d77=3094\"],\n  [\"253\u00d76=1518\", \"379\u00d78=3032\"],\n  [\"246\u00d73=738\", \"993\u00d72=1986\"],\n  [\"257\u00d79=2313\", \"208\u00d79=1872\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Old -> new three-digit x one-digit multiplication expressions\n# (each old value is unique within the document, so Find/Replace\n# per pair is unambiguous).\n$pairs = @(\n    @(\"578\u00d74=2312\", \"344\u00d77=2408\"),\n    @(\"368\u00d78=2944\", \"686\u00d76=4116\"),\n    @(\"900\u00d76=5400\", \"986\u00d75=4930\"),\n    @(\"804\u00d79=7236\", \"193\u00d72=386\"),\n    @(\"766\u00d76=4596\", \"996\u00d72=1992\"),\n    @(\"425\u00d73=1275\", \"576\u00d72=1152\"),\n    @(\"326\u00d78=2608\", \"405\u00d78=3240\"),\n    @(\"266\u00d79=2394\", \"463\u00d78=3704\"),\n    @(\"810\u00d73=2430\", \"998\u00d78=7984\"),\n    @(\"452\u00d76=2712\", \"460\u00d75=2300\"),\n    @(\"346\u00d74=1384\", \"934\u00d76=5604\"),\n    @(\"822\u00d76=4932\", \"199\u00d77=1393\"),\n    @(\"715\u00d75=3575\", \"613\u00d78=4904\"),\n    @(\"187\u00d73=561\", \"881\u00d78=7048\"),\n    @(\"178\u00d77=1246\", \"667\u00d72=1334\"),\n    @(\"305\u00d72=610\", \"391\u00d73=1173\"),\n    @(\"581\u00d74=2324\", \"238\u00d72=476\"),\n    @(\"244\u00d72=488\", \"803\u00d78=6424\"),\n    @(\"705\u00d77=4935\", \"333\u00d74=1332\"),\n    @(\"564\u00d73=1692\", \"967\u00d72=1934\"),\n    @(\"583\u00d75=2915\", \"131\u00d76=786\"),\n    @(\"357\u00d77=2499\", \"442\u00d77=3094\"),\n    @(\"253\u00d76=1518\", \"379\u00d78=3032\"),\n    @(\"246\u00d73=738\", \"993\u00d72=1986\"),\n    @(\"257\u00d79=2313\", \"208\u00d79=1872\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}\n"}
